$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Pspn"
$ws.Cells.Item(2, 3).Value = "Ret"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.166539666666667
$ws.Cells.Item(2, 8).Value = 3.499619
$ws.Cells.Item(2, 9).Value = 0.5450347805088984
$ws.Cells.Item(2, 10).Value = 0.6362259982609142
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.8639135
$ws.Cells.Item(2, 14).Value = 1.727827
$ws.Cells.Item(2, 15).Value = 0.1895490737713731
$ws.Cells.Item(2, 16).Value = 0.1380033232738433
$ws.Cells.Item(2, 17).Value = 1.007789366318833
$ws.Cells.Item(2, 18).Value = 6.046736197913
$ws.Cells.Item(2, 19).Value = 0.1033108378186453
$ws.Cells.Item(2, 20).Value = 0.08780130211322461

$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Pspn"
$ws.Cells.Item(3, 3).Value = "Ret"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.166539666666667
$ws.Cells.Item(3, 8).Value = 3.499619
$ws.Cells.Item(3, 9).Value = 0.5450347805088984
$ws.Cells.Item(3, 10).Value = 0.6362259982609142
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.384385
$ws.Cells.Item(3, 14).Value = 10.153155
$ws.Cells.Item(3, 15).Value = 0.7425593442349591
$ws.Cells.Item(3, 16).Value = 0.8109429541930055
$ws.Cells.Item(3, 17).Value = 3.948019349771666
$ws.Cells.Item(3, 18).Value = 35.53217414794499
$ws.Cells.Item(3, 19).Value = 0.4047206691999324
$ws.Cells.Item(3, 20).Value = 0.5159429905640998

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Pspn"
$ws.Cells.Item(4, 3).Value = "Ret"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.166539666666667
$ws.Cells.Item(4, 8).Value = 3.499619
$ws.Cells.Item(4, 9).Value = 0.5450347805088984
$ws.Cells.Item(4, 10).Value = 0.6362259982609142
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.2890925
$ws.Cells.Item(4, 14).Value = 0.5781849999999999
$ws.Cells.Item(4, 15).Value = 0.06342905349812297
$ws.Cells.Item(4, 16).Value = 0.04618023185601746
$ws.Cells.Item(4, 17).Value = 0.3372378685858333
$ws.Cells.Item(4, 18).Value = 2.023427211515
$ws.Cells.Item(4, 19).Value = 0.03457104025123663
$ws.Cells.Item(4, 20).Value = 0.02938106411251518

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Pspn"
$ws.Cells.Item(5, 3).Value = "Ret"
$ws.Cells.Item(5, 4).Value = "Neutrophils"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 1.166539666666667
$ws.Cells.Item(5, 8).Value = 3.499619
$ws.Cells.Item(5, 9).Value = 0.5450347805088984
$ws.Cells.Item(5, 10).Value = 0.6362259982609142
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.01599833333333333
$ws.Cells.Item(5, 14).Value = 0.047995
$ws.Cells.Item(5, 15).Value = 0.003510153811948785
$ws.Cells.Item(5, 16).Value = 0.003833410116017465
$ws.Cells.Item(5, 17).Value = 0.01866269043388889
$ws.Cells.Item(5, 18).Value = 0.167964213905
$ws.Cells.Item(5, 19).Value = 0.001913155912447979
$ws.Cells.Item(5, 20).Value = 0.002438915177806699

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Pspn"
$ws.Cells.Item(6, 3).Value = "Ret"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 1.166539666666667
$ws.Cells.Item(6, 8).Value = 3.499619
$ws.Cells.Item(6, 9).Value = 0.5450347805088984
$ws.Cells.Item(6, 10).Value = 0.6362259982609142
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.004340666666666667
$ws.Cells.Item(6, 14).Value = 0.013022
$ws.Cells.Item(6, 15).Value = 0.0009523746835961471
$ws.Cells.Item(6, 16).Value = 0.001040080561116354
$ws.Cells.Item(6, 17).Value = 0.005063559846444445
$ws.Cells.Item(6, 18).Value = 0.045572038618
$ws.Cells.Item(6, 19).Value = 0.0005190773266360577
$ws.Cells.Item(6, 20).Value = 0.0006617262932680244

$ws.Cells.Item(7, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7, 2).Value = "Pspn"
$ws.Cells.Item(7, 3).Value = "Ret"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.053445
$ws.Cells.Item(7, 8).Value = 0.160335
$ws.Cells.Item(7, 9).Value = 0.02497076154086894
$ws.Cells.Item(7, 10).Value = 0.02914868602301098
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.8639135
$ws.Cells.Item(7, 14).Value = 1.727827
$ws.Cells.Item(7, 15).Value = 0.1895490737713731
$ws.Cells.Item(7, 16).Value = 0.1380033232738433
$ws.Cells.Item(7, 17).Value = 0.0461718570075
$ws.Cells.Item(7, 18).Value = 0.277031142045
$ws.Cells.Item(7, 19).Value = 0.004733184721437532
$ws.Cells.Item(7, 20).Value = 0.004022615540241342

$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "Pspn"
$ws.Cells.Item(8, 3).Value = "Ret"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.053445
$ws.Cells.Item(8, 8).Value = 0.160335
$ws.Cells.Item(8, 9).Value = 0.02497076154086894
$ws.Cells.Item(8, 10).Value = 0.02914868602301098
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 3.384385
$ws.Cells.Item(8, 14).Value = 10.153155
$ws.Cells.Item(8, 15).Value = 0.7425593442349591
$ws.Cells.Item(8, 16).Value = 0.8109429541930055
$ws.Cells.Item(8, 17).Value = 0.180878456325
$ws.Cells.Item(8, 18).Value = 1.627906106925
$ws.Cells.Item(8, 19).Value = 0.01854227231483517
$ws.Cells.Item(8, 20).Value = 0.0236379215543449

$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "Pspn"
$ws.Cells.Item(9, 3).Value = "Ret"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.053445
$ws.Cells.Item(9, 8).Value = 0.160335
$ws.Cells.Item(9, 9).Value = 0.02497076154086894
$ws.Cells.Item(9, 10).Value = 0.02914868602301098
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.2890925
$ws.Cells.Item(9, 14).Value = 0.5781849999999999
$ws.Cells.Item(9, 15).Value = 0.06342905349812297
$ws.Cells.Item(9, 16).Value = 0.04618023185601746
$ws.Cells.Item(9, 17).Value = 0.0154505486625
$ws.Cells.Item(9, 18).Value = 0.09270329197499999
$ws.Cells.Item(9, 19).Value = 0.001583871769664647
$ws.Cells.Item(9, 20).Value = 0.001346093078840903

$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Pspn"
$ws.Cells.Item(10, 3).Value = "Ret"
$ws.Cells.Item(10, 4).Value = "Neutrophils"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.053445
$ws.Cells.Item(10, 8).Value = 0.160335
$ws.Cells.Item(10, 9).Value = 0.02497076154086894
$ws.Cells.Item(10, 10).Value = 0.02914868602301098
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.01599833333333333
$ws.Cells.Item(10, 14).Value = 0.047995
$ws.Cells.Item(10, 15).Value = 0.003510153811948785
$ws.Cells.Item(10, 16).Value = 0.003833410116017465
$ws.Cells.Item(10, 17).Value = 0.000855030925
$ws.Cells.Item(10, 18).Value = 0.007695278325000001
$ws.Cells.Item(10, 19).Value = 0.00008765121380994522
$ws.Cells.Item(10, 20).Value = 0.0001117388678692272

$ws.Cells.Item(11, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 2).Value = "Pspn"
$ws.Cells.Item(11, 3).Value = "Ret"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.053445
$ws.Cells.Item(11, 8).Value = 0.160335
$ws.Cells.Item(11, 9).Value = 0.02497076154086894
$ws.Cells.Item(11, 10).Value = 0.02914868602301098
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.004340666666666667
$ws.Cells.Item(11, 14).Value = 0.013022
$ws.Cells.Item(11, 15).Value = 0.0009523746835961471
$ws.Cells.Item(11, 16).Value = 0.001040080561116354
$ws.Cells.Item(11, 17).Value = 0.00023198693
$ws.Cells.Item(11, 18).Value = 0.00208788237
$ws.Cells.Item(11, 19).Value = 0.00002378152112163989
$ws.Cells.Item(11, 20).Value = 0.0000303169817146177

$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Pspn"
$ws.Cells.Item(12, 3).Value = "Ret"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.5
$ws.Cells.Item(12, 7).Value = 0.9203185
$ws.Cells.Item(12, 8).Value = 1.840637
$ws.Cells.Item(12, 9).Value = 0.4299944579502328
$ws.Cells.Item(12, 10).Value = 0.3346253157160749
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.8639135
$ws.Cells.Item(12, 14).Value = 1.727827
$ws.Cells.Item(12, 15).Value = 0.1895490737713731
$ws.Cells.Item(12, 16).Value = 0.1380033232738433
$ws.Cells.Item(12, 17).Value = 0.79507557644975
$ws.Cells.Item(12, 18).Value = 3.180302305799
$ws.Cells.Item(12, 19).Value = 0.08150505123129025
$ws.Cells.Item(12, 20).Value = 0.04617940562037736

$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Pspn"
$ws.Cells.Item(13, 3).Value = "Ret"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.5
$ws.Cells.Item(13, 7).Value = 0.9203185
$ws.Cells.Item(13, 8).Value = 1.840637
$ws.Cells.Item(13, 9).Value = 0.4299944579502328
$ws.Cells.Item(13, 10).Value = 0.3346253157160749
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 3.384385
$ws.Cells.Item(13, 14).Value = 10.153155
$ws.Cells.Item(13, 15).Value = 0.7425593442349591
$ws.Cells.Item(13, 16).Value = 0.8109429541930055
$ws.Cells.Item(13, 17).Value = 3.1147121266225
$ws.Cells.Item(13, 18).Value = 18.688272759735
$ws.Cells.Item(13, 19).Value = 0.3192964027201916
$ws.Cells.Item(13, 20).Value = 0.2713620420745609

$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Pspn"
$ws.Cells.Item(14, 3).Value = "Ret"
$ws.Cells.Item(14, 4).Value = "MuSCs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.5
$ws.Cells.Item(14, 7).Value = 0.9203185
$ws.Cells.Item(14, 8).Value = 1.840637
$ws.Cells.Item(14, 9).Value = 0.4299944579502328
$ws.Cells.Item(14, 10).Value = 0.3346253157160749
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.2890925
$ws.Cells.Item(14, 14).Value = 0.5781849999999999
$ws.Cells.Item(14, 15).Value = 0.06342905349812297
$ws.Cells.Item(14, 16).Value = 0.04618023185601746
$ws.Cells.Item(14, 17).Value = 0.26605717596125
$ws.Cells.Item(14, 18).Value = 1.064228703845
$ws.Cells.Item(14, 19).Value = 0.0272741414772217
$ws.Cells.Item(14, 20).Value = 0.01545307466466138

$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Pspn"
$ws.Cells.Item(15, 3).Value = "Ret"
$ws.Cells.Item(15, 4).Value = "Neutrophils"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.5
$ws.Cells.Item(15, 7).Value = 0.9203185
$ws.Cells.Item(15, 8).Value = 1.840637
$ws.Cells.Item(15, 9).Value = 0.4299944579502328
$ws.Cells.Item(15, 10).Value = 0.3346253157160749
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.01599833333333333
$ws.Cells.Item(15, 14).Value = 0.047995
$ws.Cells.Item(15, 15).Value = 0.003510153811948785
$ws.Cells.Item(15, 16).Value = 0.003833410116017465
$ws.Cells.Item(15, 17).Value = 0.01472356213583333
$ws.Cells.Item(15, 18).Value = 0.08834137281500001
$ws.Cells.Item(15, 19).Value = 0.001509346685690861
$ws.Cells.Item(15, 20).Value = 0.00128275607034154

$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Pspn"
$ws.Cells.Item(16, 3).Value = "Ret"
$ws.Cells.Item(16, 4).Value = "Resolving-Mac"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.5
$ws.Cells.Item(16, 7).Value = 0.9203185
$ws.Cells.Item(16, 8).Value = 1.840637
$ws.Cells.Item(16, 9).Value = 0.4299944579502328
$ws.Cells.Item(16, 10).Value = 0.3346253157160749
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.004340666666666667
$ws.Cells.Item(16, 14).Value = 0.013022
$ws.Cells.Item(16, 15).Value = 0.0009523746835961471
$ws.Cells.Item(16, 16).Value = 0.001040080561116354
$ws.Cells.Item(16, 17).Value = 0.003994795835666668
$ws.Cells.Item(16, 18).Value = 0.023968775014
$ws.Cells.Item(16, 19).Value = 0.0004095158358384497
$ws.Cells.Item(16, 20).Value = 0.0003480372861337124

